$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.211.10"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.672.65"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.42"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.75"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0618"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.913.90"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.70"
$ws.Range("E13").Value = "  +19.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.617"
$ws.Range("E14").Value = "  +8.40%  "
$ws.Range("D15").Value = "1.676.41"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "30.241.81"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.66"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.78"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "0.0₃0717"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  +5.03%  "
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.72"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.75"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.29"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "1.474.33"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "79.63"
$ws.Range("E38").Value = "  +15.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.589"
$ws.Range("E39").Value = "  +5.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  -7.08%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.857"
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.02"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.27"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.45"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.807.84"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.53"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("D51").Value = "0.0₆0118"
$ws.Range("E51").Value = "  +9.38%  "
